# feat: add 2022-Q4 data
#
# 1. "总计" sheet: insert a new row for 2022-Q4 above the existing 2022-Q3
#    summary row (2022-Q3 row shifts down from row 2 to row 3).
# 2. A brand new "2022-Q4" sheet (holding per-fund detail data) is inserted
#    right before the existing "2022-Q3" sheet; the old "2022-Q3" sheet
#    (with its original per-fund data) is kept, unmodified, after it.

$wb = $excel.ActiveWorkbook

# --- helper: write a value as TEXT (no residual number-format style) ----
function Set-TextCell($cell, [string]$val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# =====================================================================
# Step 1 — "总计" sheet: add the 2022-Q4 summary row, push 2022-Q3 down
# =====================================================================
$total = $wb.Worksheets.Item("总计")

# Remember the current (2022-Q3) summary row values before they're replaced
$q3Name  = $total.Cells.Item(2, 2).Value2
$q3Count = $total.Cells.Item(2, 3).Value2
$q3Value = $total.Cells.Item(2, 4).Value2

# Copy row 2's formatting down into row 3 (keeps the "A" column style)
$total.Cells.Item(2, 1).Copy()
$total.Cells.Item(3, 1).PasteSpecial(-4122)

# Row 3 becomes the (shifted-down) 2022-Q3 summary
$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(3, 2).Value = $q3Name
$total.Cells.Item(3, 3).Value = $q3Count
$total.Cells.Item(3, 4).Value = $q3Value

# Row 2 becomes the new 2022-Q4 summary
$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 4
$total.Cells.Item(2, 4).Value = 0.01

# =====================================================================
# Step 2 — duplicate the existing "2022-Q3" sheet so its original data
# survives unchanged on its own tab, then repurpose the original tab
# (same sheet identity / position) to hold the new 2022-Q4 detail data.
# =====================================================================
$q3 = $wb.Worksheets.Item("2022-Q3")

$q3.Copy($null, $q3)
$q3Copy = $wb.Worksheets.Item(3)
$q3Copy.Name = "2022-Q3-orig"

$q3.Name = "2022-Q4"

# --- overwrite the (now) "2022-Q4" sheet with the new fund detail data ---

Set-TextCell $q3.Cells.Item(2, 4) "0.27"
Set-TextCell $q3.Cells.Item(2, 5) "88.46"
Set-TextCell $q3.Cells.Item(2, 6) "4.01"
Set-TextCell $q3.Cells.Item(2, 7) "0.0108"
$q3.Cells.Item(2, 8).Value = 10

Set-TextCell $q3.Cells.Item(3, 2) "004403"
Set-TextCell $q3.Cells.Item(3, 3) "平安股息精选沪港深股票A"
Set-TextCell $q3.Cells.Item(3, 4) "0.08"
Set-TextCell $q3.Cells.Item(3, 5) "93.53"
Set-TextCell $q3.Cells.Item(3, 6) "2.84"
Set-TextCell $q3.Cells.Item(3, 7) "0.0023"
$q3.Cells.Item(3, 8).Value = 8

Set-TextCell $q3.Cells.Item(4, 2) "004404"
Set-TextCell $q3.Cells.Item(4, 3) "平安股息精选沪港深股票C"
Set-TextCell $q3.Cells.Item(4, 4) "0.03"
Set-TextCell $q3.Cells.Item(4, 5) "93.53"
Set-TextCell $q3.Cells.Item(4, 6) "2.84"
Set-TextCell $q3.Cells.Item(4, 7) "0.0009"
$q3.Cells.Item(4, 8).Value = 8

Set-TextCell $q3.Cells.Item(5, 2) "014463"
Set-TextCell $q3.Cells.Item(5, 3) "光大保德信汇佳混合C"
Set-TextCell $q3.Cells.Item(5, 4) "0.01"
Set-TextCell $q3.Cells.Item(5, 5) "88.46"
Set-TextCell $q3.Cells.Item(5, 6) "4.01"
Set-TextCell $q3.Cells.Item(5, 7) "0.0004"
$q3.Cells.Item(5, 8).Value = 10

# Header row + A-column (row index) cells: match the "总计" sheet's plain
# bold/bordered style, then (re)write values.
$total.Cells.Item(1, 2).Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

$total.Cells.Item(2, 1).Copy()
$q3.Cells.Item(2, 1).PasteSpecial(-4122)
$q3.Cells.Item(2, 1).Value = 0

$total.Cells.Item(2, 1).Copy()
$q3.Cells.Item(3, 1).PasteSpecial(-4122)
$q3.Cells.Item(3, 1).Value = 1

$total.Cells.Item(2, 1).Copy()
$q3.Cells.Item(4, 1).PasteSpecial(-4122)
$q3.Cells.Item(4, 1).Value = 2

$total.Cells.Item(2, 1).Copy()
$q3.Cells.Item(5, 1).PasteSpecial(-4122)
$q3.Cells.Item(5, 1).Value = 3

# Rename the duplicated sheet back to its original name, now trailing
# after the new "2022-Q4" sheet.
$q3Copy.Name = "2022-Q3"
